$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.756.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.320.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  +3.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.319.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.82%  "

$ws.Range("E10").Value = "  +4.08%  "

$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.58%  "

$ws.Range("E13").Value = "  +6.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "641.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.852.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.835.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.47%  "

$ws.Range("E18").Value = "  +1.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.326.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.99%  "

$ws.Range("E21").Value = "  +2.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.899"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.45%  "

$ws.Range("E24").Value = "  +1.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("E26").Value = "  +2.01%  "

$ws.Range("E27").Value = "  +5.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.32%  "

$ws.Range("E31").Value = "  +1.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "606.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.929.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.105"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.83%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.96%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.128"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "32.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0689"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.19%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.338"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.20%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0415"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.129"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.79%  "

$ws.Range("E50").Value = "  +9.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.90%  "
